# ============================================================================
# Performance.xlsx update
#  1. Populate columns E/F on sheet "20121215D" (sheet2) with new benchmark
#     data (previously-empty columns), replacing the two header strings in
#     C1/D1 with literal numeric build ids, and adding ids to E1/F1.
#  2. Fully populate the previously-empty sheet "20121215R" (sheet3) with a
#     mirrored results table (columns A-K) including formulas and the same
#     three-rule conditional formatting used on the other two sheets.
#  3. Make "20121215R" the active sheet/tab.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 20120705
$ws2 = $wb.Worksheets.Item(2)   # 20121215D
$ws3 = $wb.Worksheets.Item(3)   # 20121215R

# ----------------------------------------------------------------------
# 1. sheet2 ("20121215D") updates
# ----------------------------------------------------------------------

# Row 1 header: shrink the custom row height and replace the two text
# labels in C1/D1 with their numeric build ids; give E1/F1 ids too.
$ws2.Rows.Item(1).RowHeight = 30
$ws2.Range("C1").Value = 1205
$ws2.Range("D1").Value = 1209
$ws2.Range("E1").Value = 1210
$ws2.Range("F1").Value = 1211

# Newly-measured E/F columns for rows 2-11.
$dataEF = @(
  @(4811, 4321),
  @(4789, 4334),
  @(4796, 4326),
  @(4769, 4306),
  @(4775, 4295),
  @(4778, 4321),
  @(4775, 4314),
  @(4792, 4297),
  @(4797, 4296),
  @(4797, 4295)
)
for ($i = 0; $i -lt $dataEF.Length; $i++) {
  $r = $i + 2
  $ws2.Cells.Item($r, 5).Value = $dataEF[$i][0]
  $ws2.Cells.Item($r, 6).Value = $dataEF[$i][1]
}

# Selection / view bookkeeping for sheet2 (no longer the active tab).
$ws2.Range("A1:G16").Select()

# ----------------------------------------------------------------------
# 2. sheet3 ("20121215R") - build out the full table from scratch
# ----------------------------------------------------------------------

$ws3.Range("A1").Value = "Test"
$ws3.Range("B1").Value = 1211
$ws3.Range("C1:K1").Value = ""

$ws3.Range("A2:A11").Formula = "=ROW()-1"
$ws3.Range("A2:A11").Value = $ws3.Range("A2:A11").Value

$ws3.Range("A12").Value = "AVG"
$ws3.Range("A13").Value = "VAR"
$ws3.Range("A14").Value = "DIFF ACCEPT"
$ws3.Range("A15").Value = "Perf (Step)"
$ws3.Range("A16").Value = "Perf (Total)"

$ws3.Range("B12").Formula = "=AVERAGE(B2:B11)"
$ws3.Range("C12:G12").FormulaR1C1 = "=AVERAGE(R[-10]C:R[-1]C)"
$ws3.Range("H12:K12").FormulaR1C1 = "=AVERAGE(R[-10]C:R[-1]C)"

$ws3.Range("B13").Formula = "=_xlfn.VAR.S(B2:B11)"
$ws3.Range("C13:G13").FormulaR1C1 = "=_xlfn.VAR.S(R[-11]C:R[-2]C)"
$ws3.Range("H13:K13").FormulaR1C1 = "=_xlfn.VAR.S(R[-11]C:R[-2]C)"

$ws3.Range("B14").Value = 0
$ws3.Range("C14").Formula = "=1-_xlfn.T.TEST(B2:B11,C2:C11,2,3)"
$ws3.Range("D14").Formula = "=1-_xlfn.T.TEST(C2:C11,D2:D11,2,3)"
$ws3.Range("E14").Formula = "=1-_xlfn.T.TEST(D2:D11,E2:E11,2,3)"
$ws3.Range("F14").Formula = "=1-_xlfn.T.TEST(D2:D11,F2:F11,2,3)"
$ws3.Range("G14").Formula = "=1-_xlfn.T.TEST(F2:F11,G2:G11,2,3)"
$ws3.Range("H14").Formula = "=1-_xlfn.T.TEST(G2:G11,H2:H11,2,3)"
$ws3.Range("I14:K14").FormulaR1C1 = "=1-_xlfn.T.TEST(RC[-2]C[-2]:R[9]C[-2],RC[-1]:R[9]C[-1],2,3)"

$ws3.Range("B15").Value = 1
$ws3.Range("C15").Formula = "=B12/C12"
$ws3.Range("D15").Formula = "=C12/D12"
$ws3.Range("E15").Formula = "=D12/E12"
$ws3.Range("F15").Formula = "=D12/F12"
$ws3.Range("G15").Formula = "=F12/G12"
$ws3.Range("H15").Formula = "=G12/H12"
$ws3.Range("I15:K15").FormulaR1C1 = "=R[-3]C[-1]/R[-3]C"

$ws3.Range("B16").Value = 1
$ws3.Range("C16").Formula = "=B12/C12"
$ws3.Range("D16").Formula = "=B12/D12"
$ws3.Range("E16").Formula = "=B12/E12"
$ws3.Range("F16").Formula = "=B12/F12"
$ws3.Range("G16").Formula = "=B12/G12"
$ws3.Range("H16").Formula = "=C12/H12"
$ws3.Range("I16:K16").FormulaR1C1 = "=R[-4]C[-6]/R[-4]C"

# Styling: mirror sheet2's header (s=2) / body (s=4) cell styles.
$ws2.Range("A1:G1").Copy()
$ws3.Range("A1:G1").PasteSpecial(-4122)
$ws3.Range("H1:K1").PasteSpecial(-4122)
$ws2.Range("A2:G11").Copy()
$ws3.Range("A2:G11").PasteSpecial(-4122)
$ws3.Range("H2:K11").PasteSpecial(-4122)
$ws2.Range("A12:G16").Copy()
$ws3.Range("A12:G16").PasteSpecial(-4122)
$ws3.Range("H12:K16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Conditional formatting B15:K16 - equal / lessThan / greaterThan 1,
# matching the palette used on the other two sheets.
$rng = $ws3.Range("B15:K16")

$fcEq = $rng.FormatConditions.Add(1, 3, "1")
$fcEq.Font.Color = 0x006100
$fcEq.Interior.Color = 0xCEEFC6

$fcLt = $rng.FormatConditions.Add(1, 6, "1")
$fcLt.Font.Color = 0x0006
$fcLt.Font.Color = 0x00069C
$fcLt.Interior.Color = 0xCEC7FF

$fcGt = $rng.FormatConditions.Add(1, 5, "1")
$fcGt.Interior.Color = 0x624024

$ws3.Range("B2").Select()
$ws3.Activate()
